$wb = $excel.ActiveWorkbook

# Hunk 1: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 598.5714
$ws.Range("I103").Value = 368.82352
$ws.Range("K103").Value = 1106.47056
$ws.Range("M103").Value = -520.47056

# Hunk 2: sheet ALC
$ws.Range("H125").Value = 466
$ws.Range("I125").Value = 466
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 4194
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -1734
$ws.Range("N125").ClearContents()

# Hunk 3: sheet ALC
$ws.Range("H137").Value = 4002527.8
$ws.Range("I137").Value = 5884570
$ws.Range("J137").Value = 3187.5
$ws.Range("K137").Value = 17653710
$ws.Range("L137").Value = 9562.5
$ws.Range("M137").Value = -17651160
$ws.Range("N137").Value = -14662.5

# Hunk 4: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 14287489
$ws.Range("I102").Value = 17858824
$ws.Range("J102").Value = 2150
$ws.Range("K102").Value = 17858824
$ws.Range("L102").Value = 2150
$ws.Range("M102").Value = -17857202
$ws.Range("N102").Value = -5394

# Hunk 5: sheet ARM
$ws.Range("H132").Value = 92631.35000000001
$ws.Range("I132").Value = 63557.688
$ws.Range("J132").Value = 159085.42
$ws.Range("K132").Value = 190673.064
$ws.Range("L132").Value = 477256.26
$ws.Range("M132").Value = -188143.064
$ws.Range("N132").Value = -482316.26

# Hunk 6: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 932.8889
$ws.Range("I64").Value = 506
$ws.Range("J64").Value = 986.25
$ws.Range("K64").Value = 506
$ws.Range("L64").Value = 986.25
$ws.Range("M64").Value = -281
$ws.Range("N64").Value = -1436.25

# Hunk 7: sheet BSM
$ws.Range("H67").Value = 932.8889
$ws.Range("I67").Value = 506
$ws.Range("J67").Value = 986.25
$ws.Range("K67").Value = 506
$ws.Range("L67").Value = 986.25
$ws.Range("M67").Value = 274
$ws.Range("N67").Value = -2546.25

# Hunk 8: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1350.6346
$ws.Range("I31").Value = 946.0238000000001
$ws.Range("J31").Value = 3050
$ws.Range("K31").Value = 946.0238000000001
$ws.Range("L31").Value = 3050
$ws.Range("M31").Value = -651.0238000000001
$ws.Range("N31").Value = -3640

# Hunk 9: sheet CRP
$ws.Range("H34").Value = 1350.6346
$ws.Range("I34").Value = 946.0238000000001
$ws.Range("J34").Value = 3050
$ws.Range("K34").Value = 946.0238000000001
$ws.Range("L34").Value = 3050
$ws.Range("M34").Value = -744.0238000000001
$ws.Range("N34").Value = -3454

# Hunk 10: sheet CRP
$ws.Range("H58").Value = 27028774
$ws.Range("I58").Value = 35715428
$ws.Range("K58").Value = 35715428
$ws.Range("M58").Value = -35715225

# Hunk 11: sheet CRP
$ws.Range("H94").Value = 2958.739
$ws.Range("J94").Value = 1089.25
$ws.Range("L94").Value = 1089.25
$ws.Range("N94").Value = -1991.25

# Hunk 12: sheet CRP
$ws.Range("H136").Value = 27028774
$ws.Range("I136").Value = 35715428
$ws.Range("K136").Value = 107146284
$ws.Range("M136").Value = -107143734

# Hunk 13: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1058
$ws.Range("I34").Value = 433.33334
$ws.Range("J34").Value = 1245.4
$ws.Range("K34").Value = 1300.00002
$ws.Range("L34").Value = 3736.2
$ws.Range("M34").Value = -1216.00002
$ws.Range("N34").Value = -3904.2

# Hunk 14: sheet CUL
$ws.Range("H39").Value = 1667.6364
$ws.Range("J39").Value = 1667.6364
$ws.Range("L39").Value = 5002.9092
$ws.Range("N39").Value = -5590.9092

# Hunk 15: sheet CUL
$ws.Range("H55").Value = 2280
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2280
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 6840
$ws.Range("N55").Value = -7194
$ws.Range("M55").ClearContents()

# Hunk 16: sheet CUL
$ws.Range("H131").Value = 2991.0566
$ws.Range("I131").Value = 569.8889
$ws.Range("J131").Value = 3486.2954
$ws.Range("K131").Value = 1709.6667
$ws.Range("L131").Value = 10458.8862
$ws.Range("M131").Value = 3330.3333
$ws.Range("N131").Value = -20538.8862

# Hunk 17: sheet CUL
$ws.Range("H132").Value = 965
$ws.Range("J132").Value = 994.44446
$ws.Range("L132").Value = 8950.00014
$ws.Range("N132").Value = -14010.00014

# Hunk 18: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1501.25
$ws.Range("I97").Value = 1646.9231
$ws.Range("J97").Value = 870
$ws.Range("K97").Value = 1646.9231
$ws.Range("L97").Value = 870
$ws.Range("M97").Value = -1150.9231
$ws.Range("N97").Value = -1862

# Hunk 19: sheet GSM
$ws.Range("H102").Value = 3368.5881
$ws.Range("I102").Value = 2804.8572
$ws.Range("J102").Value = 5999.3335
$ws.Range("K102").Value = 2804.8572
$ws.Range("L102").Value = 5999.3335
$ws.Range("M102").Value = -1182.8572
$ws.Range("N102").Value = -9243.333500000001

# Hunk 20: sheet GSM
$ws.Range("H122").Value = 2092.0667
$ws.Range("I122").Value = 1615.0834
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4845.2502
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2395.2502
$ws.Range("N122").Value = -16900

# Hunk 21: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 119597.234
$ws.Range("I136").Value = 68743.53
$ws.Range("J136").Value = 501000
$ws.Range("K136").Value = 206230.59
$ws.Range("L136").Value = 1503000
$ws.Range("M136").Value = -203680.59
$ws.Range("N136").Value = -1508100

# Hunk 22: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 12394.4
$ws.Range("J54").Value = 12394.4
$ws.Range("L54").Value = 12394.4
$ws.Range("N54").Value = -13434.4

# Hunk 23: sheet WVR
$ws.Range("H122").Value = 2318.5715
$ws.Range("I122").Value = 1085.0769
$ws.Range("J122").Value = 4323
$ws.Range("K122").Value = 3255.2307
$ws.Range("L122").Value = 12969
$ws.Range("M122").Value = -805.2307000000001
$ws.Range("N122").Value = -17869

# Hunk 24: sheet WVR
$ws.Range("H126").Value = 2348.2144
$ws.Range("I126").Value = 1797.4
$ws.Range("J126").Value = 3725.25
$ws.Range("K126").Value = 5392.200000000001
$ws.Range("L126").Value = 11175.75
$ws.Range("M126").Value = -2922.200000000001
$ws.Range("N126").Value = -16115.75

# Hunk 25: sheet WVR
$ws.Range("H132").Value = 51141.574
$ws.Range("I132").Value = 32045.781
$ws.Range("J132").Value = 127524.75
$ws.Range("K132").Value = 96137.34299999999
$ws.Range("L132").Value = 382574.25
$ws.Range("M132").Value = -93607.34299999999
$ws.Range("N132").Value = -387634.25

# Hunk 26: sheet WVR
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Hunk 27: sheet WVR
$ws.Range("H136").Value = 43910.15
$ws.Range("I136").Value = 24780.404
$ws.Range("J136").Value = 204600
$ws.Range("K136").Value = 74341.212
$ws.Range("L136").Value = 613800
$ws.Range("M136").Value = -71791.212
$ws.Range("N136").Value = -618900
